# Update the NeXus filewriter json-config table to match the updated
# filewriter format: rename "writer_module" -> "module" and add two new
# columns ("dtype" and "value_units") describing each source's data type
# and unit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---------------------------------------------------
# Column E header: "writer_module" -> "module"
$ws.Range("E1").Value = "module"
# New columns: F = dtype, G = value_units
$ws.Range("F1").Value = "dtype"
$ws.Range("G1").Value = "value_units"

# --- Data rows: new "dtype" / "value_units" columns -----------------------
$dataRange = $ws.Range("F2:G4")
$dataRange.HorizontalAlignment = -4108

$ws.Range("F2").Value = "string"
$ws.Range("G2").Value = "mm"

$ws.Range("F3").Value = "string"
$ws.Range("G3").Value = "mm"

$ws.Range("F4").Value = "string"
$ws.Range("G4").Value = "mm"

# --- Header styling: red font, centered -----------------------------------
$headerRange = $ws.Range("A1:G1")
$headerRange.HorizontalAlignment = -4108
$headerRange.Font.Color = 255

# --- View state: zoom + selection ------------------------------------------
$ws.Application.ActiveWindow.Zoom = 234
$ws.Range("F3").Select() | Out-Null
